$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = 'ECs'
$data[0,1] = 'Fgf2'
$data[0,2] = 'Fgfr3'
$data[0,3] = 'ECs'
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.472738
$data[0,7] = 1.418214
$data[0,8] = 0.0327564895931267
$data[0,9] = 0.03397138804734427
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 2.330840333333333
$data[0,13] = 6.992521
$data[0,14] = 0.6715345129768794
$data[0,15] = 0.7003397275969581
$data[0,16] = 1.101876797499333
$data[0,17] = 9.916891177494001
$data[0,18] = 0.02199711328575256
$data[0,19] = 0.02379151265116764
$data[1,0] = 'ECs'
$data[1,1] = 'Fgf2'
$data[1,2] = 'Fgfr3'
$data[1,3] = 'FAPs'
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.472738
$data[1,7] = 1.418214
$data[1,8] = 0.0327564895931267
$data[1,9] = 0.03397138804734427
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.657666
$data[1,13] = 1.972998
$data[1,14] = 0.1894790521235985
$data[1,15] = 0.1976066831789769
$data[1,16] = 0.310903709508
$data[1,17] = 2.798133385572
$data[1,18] = 0.006206668599002166
$data[1,19] = 0.006712973315021644
$data[2,0] = 'ECs'
$data[2,1] = 'Fgf2'
$data[2,2] = 'Fgfr3'
$data[2,3] = 'MuSCs'
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.472738
$data[2,7] = 1.418214
$data[2,8] = 0.0327564895931267
$data[2,9] = 0.03397138804734427
$data[2,10] = 2
$data[2,11] = 1
$data[2,12] = 0.42828
$data[2,13] = 0.85656
$data[2,14] = 0.123391035029171
$data[2,15] = 0.08578923067523865
$data[2,16] = 0.20246423064
$data[2,17] = 1.21478538384
$data[2,18] = 0.004041857154818172
$data[2,19] = 0.002914379245551663
$data[3,0] = 'ECs'
$data[3,1] = 'Fgf2'
$data[3,2] = 'Fgfr3'
$data[3,3] = 'Resolving-Mac'
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.472738
$data[3,7] = 1.418214
$data[3,8] = 0.0327564895931267
$data[3,9] = 0.03397138804734427
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.05413033333333334
$data[3,13] = 0.162391
$data[3,14] = 0.01559539987035126
$data[3,15] = 0.01626435854882633
$data[3,16] = 0.02558946551933334
$data[3,17] = 0.230305189674
$data[3,18] = 0.0005108505535538104
$data[3,19] = 0.0005525228356033203
$data[4,0] = 'FAPs'
$data[4,1] = 'Fgf2'
$data[4,2] = 'Fgfr3'
$data[4,3] = 'ECs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 11.86561333333333
$data[4,7] = 35.59684
$data[4,8] = 0.822180234441485
$data[4,9] = 0.8526739017519405
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.330840333333333
$data[4,13] = 6.992521
$data[4,14] = 0.6715345129768794
$data[4,15] = 0.7003397275969581
$data[4,16] = 27.65685013707111
$data[4,17] = 248.91165123364
$data[4,18] = 0.5521224033148792
$data[4,19] = 0.5971614080819894
$data[5,0] = 'FAPs'
$data[5,1] = 'Fgf2'
$data[5,2] = 'Fgfr3'
$data[5,3] = 'FAPs'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 11.86561333333333
$data[5,7] = 35.59684
$data[5,8] = 0.822180234441485
$data[5,9] = 0.8526739017519405
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.657666
$data[5,13] = 1.972998
$data[5,14] = 0.1894790521235985
$data[5,15] = 0.1976066831789769
$data[5,16] = 7.803610458480001
$data[5,17] = 70.23249412632001
$data[5,18] = 0.1557859314967306
$data[5,19] = 0.1684940615584778
$data[6,0] = 'FAPs'
$data[6,1] = 'Fgf2'
$data[6,2] = 'Fgfr3'
$data[6,3] = 'MuSCs'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 11.86561333333333
$data[6,7] = 35.59684
$data[6,8] = 0.822180234441485
$data[6,9] = 0.8526739017519405
$data[6,10] = 2
$data[6,11] = 1
$data[6,12] = 0.42828
$data[6,13] = 0.85656
$data[6,14] = 0.123391035029171
$data[6,15] = 0.08578923067523865
$data[6,16] = 5.0818048784
$data[6,17] = 30.4908292704
$data[6,18] = 0.1014496701082613
$data[6,19] = 0.073150238048153
$data[7,0] = 'FAPs'
$data[7,1] = 'Fgf2'
$data[7,2] = 'Fgfr3'
$data[7,3] = 'Resolving-Mac'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 11.86561333333333
$data[7,7] = 35.59684
$data[7,8] = 0.822180234441485
$data[7,9] = 0.8526739017519405
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.05413033333333334
$data[7,13] = 0.162391
$data[7,14] = 0.01559539987035126
$data[7,15] = 0.01626435854882633
$data[7,16] = 0.6422896049377779
$data[7,17] = 5.78060644444
$data[7,18] = 0.0128222295216141
$data[7,19] = 0.01386819406332027
$data[8,0] = 'Inflammatory-Mac'
$data[8,1] = 'Fgf2'
$data[8,2] = 'Fgfr3'
$data[8,3] = 'ECs'
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.37892
$data[8,7] = 1.13676
$data[8,8] = 0.02625574638939025
$data[8,9] = 0.02722954016579943
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 2.330840333333333
$data[8,13] = 6.992521
$data[8,14] = 0.6715345129768794
$data[8,15] = 0.7003397275969581
$data[8,16] = 0.8832020191066665
$data[8,17] = 7.94881817196
$data[8,18] = 0.01763163986444364
$data[8,19] = 0.0190699287423064
$data[9,0] = 'Inflammatory-Mac'
$data[9,1] = 'Fgf2'
$data[9,2] = 'Fgfr3'
$data[9,3] = 'FAPs'
$data[9,4] = 2
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.37892
$data[9,7] = 1.13676
$data[9,8] = 0.02625574638939025
$data[9,9] = 0.02722954016579943
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.657666
$data[9,13] = 1.972998
$data[9,14] = 0.1894790521235985
$data[9,15] = 0.1976066831789769
$data[9,16] = 0.24920280072
$data[9,17] = 2.24282520648
$data[9,18] = 0.004974913938659258
$data[9,19] = 0.005380739116652355
$data[10,0] = 'Inflammatory-Mac'
$data[10,1] = 'Fgf2'
$data[10,2] = 'Fgfr3'
$data[10,3] = 'MuSCs'
$data[10,4] = 2
$data[10,5] = 0.6666666666666666
$data[10,6] = 0.37892
$data[10,7] = 1.13676
$data[10,8] = 0.02625574638939025
$data[10,9] = 0.02722954016579943
$data[10,10] = 2
$data[10,11] = 1
$data[10,12] = 0.42828
$data[10,13] = 0.85656
$data[10,14] = 0.123391035029171
$data[10,15] = 0.08578923067523865
$data[10,16] = 0.1622838576
$data[10,17] = 0.9737031455999999
$data[10,18] = 0.003239723722450282
$data[10,19] = 0.002336001302464444
$data[11,0] = 'Inflammatory-Mac'
$data[11,1] = 'Fgf2'
$data[11,2] = 'Fgfr3'
$data[11,3] = 'Resolving-Mac'
$data[11,4] = 2
$data[11,5] = 0.6666666666666666
$data[11,6] = 0.37892
$data[11,7] = 1.13676
$data[11,8] = 0.02625574638939025
$data[11,9] = 0.02722954016579943
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.05413033333333334
$data[11,13] = 0.162391
$data[11,14] = 0.01559539987035126
$data[11,15] = 0.01626435854882633
$data[11,16] = 0.02051106590666667
$data[11,17] = 0.18459959316
$data[11,18] = 0.0004094688638370721
$data[11,19] = 0.0004428710043762298
$data[12,0] = 'MuSCs'
$data[12,1] = 'Fgf2'
$data[12,2] = 'Fgfr3'
$data[12,3] = 'ECs'
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 1.548357
$data[12,7] = 3.096714
$data[12,8] = 0.1072872076222874
$data[12,9] = 0.0741775733180209
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 2.330840333333333
$data[12,13] = 6.992521
$data[12,14] = 0.6715345129768794
$data[12,15] = 0.7003397275969581
$data[12,16] = 3.608972945999
$data[12,17] = 21.653837675994
$data[12,18] = 0.07204706271928209
$data[12,19] = 0.05194950149134614
$data[13,0] = 'MuSCs'
$data[13,1] = 'Fgf2'
$data[13,2] = 'Fgfr3'
$data[13,3] = 'FAPs'
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 1.548357
$data[13,7] = 3.096714
$data[13,8] = 0.1072872076222874
$data[13,9] = 0.0741775733180209
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.657666
$data[13,13] = 1.972998
$data[13,14] = 0.1894790521235985
$data[13,15] = 0.1976066831789769
$data[13,16] = 1.018301754762
$data[13,17] = 6.109810528572001
$data[13,18] = 0.02032867840525872
$data[13,19] = 0.01465798422963949
$data[14,0] = 'MuSCs'
$data[14,1] = 'Fgf2'
$data[14,2] = 'Fgfr3'
$data[14,3] = 'MuSCs'
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 1.548357
$data[14,7] = 3.096714
$data[14,8] = 0.1072872076222874
$data[14,9] = 0.0741775733180209
$data[14,10] = 2
$data[14,11] = 1
$data[14,12] = 0.42828
$data[14,13] = 0.85656
$data[14,14] = 0.123391035029171
$data[14,15] = 0.08578923067523865
$data[14,16] = 0.66313033596
$data[14,17] = 2.65252134384
$data[14,18] = 0.0132382795939036
$data[14,19] = 0.006363636948309123
$data[15,0] = 'MuSCs'
$data[15,1] = 'Fgf2'
$data[15,2] = 'Fgfr3'
$data[15,3] = 'Resolving-Mac'
$data[15,4] = 2
$data[15,5] = 1
$data[15,6] = 1.548357
$data[15,7] = 3.096714
$data[15,8] = 0.1072872076222874
$data[15,9] = 0.0741775733180209
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.05413033333333334
$data[15,13] = 0.162391
$data[15,14] = 0.01559539987035126
$data[15,15] = 0.01626435854882633
$data[15,16] = 0.08381308052900001
$data[15,17] = 0.5028784831740001
$data[15,18] = 0.001673186903842968
$data[15,19] = 0.001206450648726145
$data[16,0] = 'Resolving-Mac'
$data[16,1] = 'Fgf2'
$data[16,2] = 'Fgfr3'
$data[16,3] = 'ECs'
$data[16,4] = 1
$data[16,5] = 0.3333333333333333
$data[16,6] = 0.16626
$data[16,7] = 0.49878
$data[16,8] = 0.01152032195371061
$data[16,9] = 0.01194759671689489
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 2.330840333333333
$data[16,13] = 6.992521
$data[16,14] = 0.6715345129768794
$data[16,15] = 0.7003397275969581
$data[16,16] = 0.3875255138199999
$data[16,17] = 3.48772962438
$data[16,18] = 0.007736293792521904
$data[16,19] = 0.00836737663014848
$data[17,0] = 'Resolving-Mac'
$data[17,1] = 'Fgf2'
$data[17,2] = 'Fgfr3'
$data[17,3] = 'FAPs'
$data[17,4] = 1
$data[17,5] = 0.3333333333333333
$data[17,6] = 0.16626
$data[17,7] = 0.49878
$data[17,8] = 0.01152032195371061
$data[17,9] = 0.01194759671689489
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 0.657666
$data[17,13] = 1.972998
$data[17,14] = 0.1894790521235985
$data[17,15] = 0.1976066831789769
$data[17,16] = 0.10934354916
$data[17,17] = 0.98409194244
$data[17,18] = 0.002182859683947768
$data[17,19] = 0.002360924959185634
$data[18,0] = 'Resolving-Mac'
$data[18,1] = 'Fgf2'
$data[18,2] = 'Fgfr3'
$data[18,3] = 'MuSCs'
$data[18,4] = 1
$data[18,5] = 0.3333333333333333
$data[18,6] = 0.16626
$data[18,7] = 0.49878
$data[18,8] = 0.01152032195371061
$data[18,9] = 0.01194759671689489
$data[18,10] = 2
$data[18,11] = 1
$data[18,12] = 0.42828
$data[18,13] = 0.85656
$data[18,14] = 0.123391035029171
$data[18,15] = 0.08578923067523865
$data[18,16] = 0.07120583279999999
$data[18,17] = 0.4272349968
$data[18,18] = 0.001421504449737633
$data[18,19] = 0.00102497513076042
$data[19,0] = 'Resolving-Mac'
$data[19,1] = 'Fgf2'
$data[19,2] = 'Fgfr3'
$data[19,3] = 'Resolving-Mac'
$data[19,4] = 1
$data[19,5] = 0.3333333333333333
$data[19,6] = 0.16626
$data[19,7] = 0.49878
$data[19,8] = 0.01152032195371061
$data[19,9] = 0.01194759671689489
$data[19,10] = 1
$data[19,11] = 0.3333333333333333
$data[19,12] = 0.05413033333333334
$data[19,13] = 0.162391
$data[19,14] = 0.01559539987035126
$data[19,15] = 0.01626435854882633
$data[19,16] = 0.008999709220000001
$data[19,17] = 0.08099738298
$data[19,18] = 0.0001796640275033031
$data[19,19] = 0.0001943199968003588

$ws.Range("A2:T21").Value = $data